$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the existing "Mut Score" header (H4) to "Mut Score Botium" ---
$ws.Range("H4").Value = "Mut Score Botium"

# --- Clear the old, now-superseded "#tests = 52" header cell (content only, keep its format) ---
$ws.Range("G1").ClearContents()

# --- Normalize "Killed" -> "KILLED" wording for consistency ---
$ws.Range("E2").Value = "KILLED"

# --- New summary labels in columns H:I ---
$ws.Range("H3").Value = "Botium # tests = 52"
$ws.Range("H2").Value = "mutantes = "

# Column H needs to be wide enough for the new labels.
$ws.Columns.Item(8).ColumnWidth = 22.6

# I2: total mutants analysed so far (killed + survived), styled like the
# highlighted "Done" cells (fill, no border).
$ws.Range("A2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("I2").Borders.LineStyle = -4142
$ws.Range("I2").Formula = '=(COUNTIF(E2:E94, "Killed")+COUNTIF(E2:E94, "survived"))'


# --- Fill in "Done" status + mutant verdict for rows that were pending ---
$ws.Range("B33").Value = "Done"
$ws.Range("E33").Value = "survived"

$ws.Range("B36").Value = "Done"
$ws.Range("E36").Value = "KILLED"

$ws.Range("B39").Value = "Done"
$ws.Range("E39").Value = "KILLED"

$ws.Range("B51").Value = "Done"
$ws.Range("E51").Value = "survived"

$ws.Range("B88").Value = "Done"
$ws.Range("E88").Value = "survived"
# Row 88's F cell had an inconsistent format (no fill) vs the rest of the row; fix it to match.
$ws.Range("B88").Copy()
$ws.Range("F88").PasteSpecial(-4122)
$ws.Range("F88").Value = ""

$ws.Range("I4").Select()

$wb.Save()
